# Update "想去人数" (want-to-go count) figures in column F, regenerated output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 667
$ws1.Range("F3").Value = 46
$ws1.Range("F4").Value = 1985
$ws1.Range("F5").Value = 5751
$ws1.Range("F6").Value = 1608
$ws1.Range("F7").Value = 170
$ws1.Range("F8").Value = 3252
$ws1.Range("F11").Value = 1366
$ws1.Range("F12").Value = 4538
$ws1.Range("F13").Value = 1083
$ws1.Range("F14").Value = 1710
$ws1.Range("F15").Value = 2607
$ws1.Range("F17").Value = 48
$ws1.Range("F18").Value = 52
$ws1.Range("F19").Value = 179
$ws1.Range("F21").Value = 1023
$ws1.Range("F22").Value = 307
$ws1.Range("F29").Value = 1119
$ws1.Range("F31").Value = 86
$ws1.Range("F32").Value = 203
$ws1.Range("F33").Value = 380
$ws1.Range("F34").Value = 839
$ws1.Range("F36").Value = 1742
$ws1.Range("F37").Value = 2253
$ws1.Range("F38").Value = 1050
$ws1.Range("F40").Value = 272
$ws1.Range("F42").Value = 372
$ws1.Range("F44").Value = 671
$ws1.Range("F46").Value = 442
$ws1.Range("F47").Value = 406
$ws1.Range("F48").Value = 231

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 13
$ws2.Range("F8").Value = 1

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 667
$ws4.Range("F4").Value = 46
$ws4.Range("F5").Value = 1985
$ws4.Range("F6").Value = 5751
$ws4.Range("F7").Value = 1608
$ws4.Range("F8").Value = 170
$ws4.Range("F9").Value = 3252
$ws4.Range("F11").Value = 1366
$ws4.Range("F12").Value = 4538
$ws4.Range("F13").Value = 1083
$ws4.Range("F14").Value = 1710
$ws4.Range("F17").Value = 48
$ws4.Range("F19").Value = 52
$ws4.Range("F20").Value = 179
$ws4.Range("F23").Value = 1023
$ws4.Range("F24").Value = 307
$ws4.Range("F29").Value = 1119
$ws4.Range("F31").Value = 86
$ws4.Range("F32").Value = 203
$ws4.Range("F33").Value = 840
$ws4.Range("F34").Value = 1742
$ws4.Range("F35").Value = 2253
$ws4.Range("F36").Value = 1050
$ws4.Range("F40").Value = 272
$ws4.Range("F42").Value = 372
$ws4.Range("F43").Value = 671
$ws4.Range("F44").Value = 442
$ws4.Range("F45").Value = 406
$ws4.Range("F46").Value = 231
